# Commit: "From 0.1 to 1.0 version"
#
# Net semantic change in the sheet:
#  1) The "Version: " value (cell D2) goes from "0.1" to "1.0".
#  2) Step 1 of every test case ("Usuario do Sistema inicia a tela de login
#     atraves da opcao de login no canto superior direito") capitalizes the
#     word "login" -> "Login" in "...opcao de login no canto...". Because
#     that sentence is a single shared string reused by all six test cases
#     (TC1-TC6, rows 10/20/30/40/50/60, column B), editing it once updates
#     every test case at once - matching the diff, which only touches the
#     shared-strings table plus the index shuffles that fall out of it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# 1) Bump the version number.
$ws.Range("D2").Value = "1.0"

# 2) Capitalize "Login" in the shared "step 1" sentence for every test case.
$oldText = "Usuario do Sistema inicia a tela de login atraves da opcao de login no canto superior direito"
$newText = "Usuario do Sistema inicia a tela de login atraves da opcao de Login no canto superior direito"

foreach ($rowNum in 10, 20, 30, 40, 50, 60) {
    $cell = $ws.Range("B$rowNum")
    if ($cell.Value -eq $oldText) {
        $cell.Value = $newText
    }
}
